$d = $word.ActiveDocument

# The edit touches the "Leverrier" section (removing proofErr spell-check
# markers and merging the runs they split) and restructures the paragraphs
# between that section and "Exploitation" (new blank paragraphs, the
# "Presentation du code" heading takes over the old _GoBack marker's slot,
# bookmark ids shift down by one, and a fresh _GoBack bookmark paragraph is
# inserted right after "Exploitation").
#
# Paragraph 21 is "La methode des puissances..." and paragraph 33 is
# "Exploitation" in the original document - together they span exactly the
# block touched by the diff. Replacing that whole block in one InsertXML
# call lets us control proofErr markers, run splits, bookmark ids and
# paragraph counts precisely.

$startPara = $d.Paragraphs.Item(21)
$endPara = $d.Paragraphs.Item(33)
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00361BE3" w:rsidRPr="00361BE3" w:rsidRDefault="00361BE3" w:rsidP="00361BE3"><w:r><w:t xml:space="preserve">La </w:t></w:r><w:r w:rsidR="00460E65"><w:t>méthode</w:t></w:r><w:r><w:t xml:space="preserve"> des puissances consiste à définir une Suite en fonction de </w:t></w:r><w:r w:rsidR="00460E65"><w:t>(</w:t></w:r><w:r><w:t>X</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>n</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r w:rsidR="00460E65"><w:t>= A*</w:t></w:r><w:r w:rsidR="00460E65" w:rsidRPr="00460E65"><w:t>X</w:t></w:r><w:r w:rsidR="00460E65"><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>n-1</w:t></w:r><w:r w:rsidR="00460E65"><w:t xml:space="preserve"> qui</w:t></w:r><w:r><w:t xml:space="preserve"> converge </w:t></w:r><w:r w:rsidR="00460E65"><w:t xml:space="preserve">vers le vecteur </w:t></w:r><w:r w:rsidR="00D83BF9"><w:t>propre de la matrice A</w:t></w:r><w:r w:rsidR="00E61A46"><w:t>.</w:t></w:r><w:r w:rsidR="00457F5D"><w:t xml:space="preserve">  À partir de ce vecteur propre on peut calculer la plus grande valeur propre de la matrice.</w:t></w:r></w:p><w:p w:rsidR="002046A8" w:rsidRDefault="00AB3850" w:rsidP="002046A8"><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:bookmarkStart w:id="2" w:name="_Toc469666494"/><w:r><w:t>Le</w:t></w:r><w:r w:rsidR="00E61A46"><w:t>verrier</w:t></w:r><w:bookmarkEnd w:id="2"/></w:p><w:p w:rsidR="00365F9C" w:rsidRDefault="0023422E" w:rsidP="0023422E"><w:r><w:t>La méthode de Leverrier sert à obtenir un polynôme caractéristique à partir d’une matrice carrée. Les valeurs pour lesquelles ce polynôme s’annule sont les valeurs propres de la matrice. Pour avoir les coeffici</w:t></w:r><w:r><w:t xml:space="preserve">ents de ce polynôme, Leverrier a créé une formule qui se sert des traces de la matrice de la puissance 1 à n.  Le Polynôme caractéristique est de la forme : </w:t></w:r></w:p><w:p w:rsidR="00365F9C" w:rsidRPr="0023422E" w:rsidRDefault="00365F9C" w:rsidP="0023422E"><w:r><w:t>P(X) = | A-XIn |= an + an-1X+an-2X^2+…+a0X^n</w:t></w:r></w:p><w:p w:rsidR="00B75122" w:rsidRPr="00B75122" w:rsidRDefault="00B75122" w:rsidP="00B75122"/><w:p/><w:p w:rsidR="000F4103" w:rsidRDefault="000F4103" w:rsidP="000F4103"><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:bookmarkStart w:id="3" w:name="_Toc469666495"/><w:r><w:t>Présentation du code</w:t></w:r><w:bookmarkEnd w:id="3"/></w:p><w:p/><w:p w:rsidR="008B3CFB" w:rsidRDefault="008B3CFB" w:rsidP="008B3CFB"><w:pPr><w:pStyle w:val="Titre2"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:bookmarkStart w:id="4" w:name="_Toc469666496"/><w:r><w:t>Choix de programmation</w:t></w:r><w:bookmarkEnd w:id="4"/></w:p><w:p/><w:p w:rsidR="00D1352F" w:rsidRDefault="009C4889" w:rsidP="009C4889"><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:bookmarkStart w:id="5" w:name="_Toc469666497"/><w:r><w:t>Organisation du code</w:t></w:r><w:bookmarkEnd w:id="5"/></w:p><w:p w:rsidR="009C4889" w:rsidRDefault="00287BCE" w:rsidP="009C4889"><w:r><w:br/></w:r></w:p><w:p w:rsidR="00D6792C" w:rsidRDefault="00D6792C" w:rsidP="006E3238"><w:pPr><w:pStyle w:val="Titre3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="360"/></w:pPr></w:p><w:p w:rsidR="00A62254" w:rsidRDefault="00DB0573" w:rsidP="00A62254"><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:bookmarkStart w:id="6" w:name="_Toc469666498"/><w:r><w:t>Améliorations</w:t></w:r><w:r w:rsidR="009F5606"><w:t> :</w:t></w:r><w:bookmarkEnd w:id="6"/><w:r w:rsidR="009F5606"><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w:rsidR="00C23FF8" w:rsidRDefault="004F7409" w:rsidP="003D667C"><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:bookmarkStart w:id="7" w:name="_Toc469666499"/><w:r><w:t>E</w:t></w:r><w:r w:rsidR="003D667C"><w:t>xploitation</w:t></w:r><w:bookmarkEnd w:id="7"/></w:p><w:p><w:bookmarkStart w:id="8" w:name="_GoBack"/><w:bookmarkEnd w:id="8"/></w:p>
'@

$r.InsertXML($xml)

Write-Output "done"
